# Apply cryptos list update (GitHub Actions refresh) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '36.478.33'
$ws.Range('E2').Value = '  -2.63%  '

$ws.Range('D3').Value = '1.952.81'
$ws.Range('E3').Value = '  -3.47%  '

$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.88'
$ws.Range('E5').Value = '  -10.81%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.593'
$ws.Range('E6').Value = '  -4.24%  '

$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.42'
$ws.Range('E8').Value = '  -7.43%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.362'
$ws.Range('E9').Value = '  -5.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.82'
$ws.Range('E10').Value = '  -0.72%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0726'
$ws.Range('E11').Value = '  -7.18%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0968'
$ws.Range('E12').Value = '  -4.56%  '

$ws.Range('D13').Value = '2.245.08'
$ws.Range('E13').Value = '  -3.21%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '13.58'
$ws.Range('E14').Value = '  -6.43%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.29'
$ws.Range('E15').Value = '  -8.74%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.731'
$ws.Range('E16').Value = '  -10.35%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.94'
$ws.Range('E17').Value = '  -7.59%  '

$ws.Range('D18').Value = '1.958.28'
$ws.Range('E18').Value = '  -3.08%  '

$ws.Range('D19').Value = '36.467.25'
$ws.Range('E19').Value = '  -2.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '66.79'
$ws.Range('E20').Value = '  -3.80%  '

$ws.Range('D21').Value = '0.0₃0781'
$ws.Range('E21').Value = '  -7.92%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.88'
$ws.Range('E22').Value = '  -5.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '217.73'
$ws.Range('E23').Value = '  -4.64%  '

$ws.Range('E24').Value = '  +0.02%  '

$ws.Range('E25').Value = '  -0.81%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  -12.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.51'
$ws.Range('E27').Value = '  -2.72%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.35'
$ws.Range('E28').Value = '  -7.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.54'
$ws.Range('E29').Value = '  -6.61%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.120'
$ws.Range('E30').Value = '  -8.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.27'
$ws.Range('E31').Value = '  -7.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  -4.21%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.23'
$ws.Range('E33').Value = '  -10.06%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0591'
$ws.Range('E34').Value = '  -10.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.11'
$ws.Range('E35').Value = '  -10.02%  '

$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.31%  '

$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.23'
$ws.Range('E37').Value = '  -8.32%  '

$ws.Range('E38').Value = '  -2.20%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.13'
$ws.Range('E39').Value = '  -7.91%  '

$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.06'
$ws.Range('E40').Value = '  +0.93%  '

$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.04'
$ws.Range('E41').Value = '  -5.44%  '

$ws.Range('D42').Value = '1.389.99'
$ws.Range('E42').Value = '  -0.45%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0198'
$ws.Range('E43').Value = '  -8.04%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0863'
$ws.Range('E44').Value = '  -10.39%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.06'
$ws.Range('E45').Value = '  -10.73%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '85.11'
$ws.Range('E46').Value = '  -6.17%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.968'
$ws.Range('E47').Value = '  -6.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '14.55'
$ws.Range('E48').Value = '  -8.95%  '

$ws.Range('E49').Value = '  -0.39%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.61'
$ws.Range('E50').Value = '  -9.68%  '

$ws.Range('D51').Value = '2.138.14'
$ws.Range('E51').Value = '  -3.33%  '
